$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing sheet "journalVoucherDetails" ---
$ws1.Range("B2").Value = "27/12/2016"
$ws1.Range("B3").Value = "27/12/2016"
$ws1.Range("A4").Value = "voucherDateJune"

# restore the selection on the (now inactive) first sheet
$ws1.Range("A14").Select() | Out-Null

# --- Add the new sheet "financialBankDetails" right after it ---
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "financialBankDetails"

$newSheet.Range("A1").Value = "dataName"
$newSheet.Range("B1").Value = "bankName"
$newSheet.Range("C1").Value = "accountNumber"
$newSheet.Range("A2").Value = "SBI"
$newSheet.Range("B2").Value = "ANDHRA BANK-Andhra Bank RTC Busstand"
$newSheet.Range("C2").Value = "4502110--110710100009664--ANDHRA BANK"

# accountNumber column is text-formatted, like the numeric-looking account codes on sheet 1
$newSheet.Columns.Item(3).NumberFormat = "@"

# make the new sheet the active / selected tab
$newSheet.Activate() | Out-Null
$newSheet.Range("C2").Select() | Out-Null
